# MOM_03_03_2023.docx weekly update edit
# Applies the "mom for this week" commit: rewrites bullet content across the
# "Specific Activities", "Specific Output from work", "On Target",
# "Challenges/Disagreements" and "Planned Activities for coming week" sections.

$d = $word.ActiveDocument

function Find-ParaIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# Rewrites the text of paragraph $paraIndex using a list of runs. Each run is
# a hashtable like @{Text="abc"; Bold=$true}. The paragraph mark itself is
# left untouched (only the paragraph's own text content is replaced).
#
# NOTE: we intentionally delete the old text (excluding the trailing
# paragraph mark) and then insert the new text, rather than assigning
# $r.Text directly — when a Range spans more than one pre-existing run,
# assigning .Text only overwrites the first run and silently leaves the
# remaining old runs behind.
function Set-ParaRuns($doc, $paraIndex, $segments) {
    Set-ParaRunsKeepTail $doc $paraIndex $segments 0
}

# Same as Set-ParaRuns, but preserves the last $tailChars characters of the
# paragraph's existing content untouched (used where the diff shows a
# trailing run, e.g. a lone "." or ":", surviving unchanged).
function Set-ParaRunsKeepTail($doc, $paraIndex, $segments, $tailChars) {
    $p = $doc.Paragraphs($paraIndex)
    $start = $p.Range.Start
    $end = $p.Range.End
    $replaceEnd = $end - 1 - $tailChars

    if ($replaceEnd -gt $start) {
        $doc.Range($start, $replaceEnd).Delete()
    }

    $pos = $start
    foreach ($seg in $segments) {
        $len = $seg.Text.Length
        if ($len -gt 0) {
            $ins = $doc.Range($pos, $pos)
            $ins.InsertBefore($seg.Text)
            if ($seg.Bold) {
                $sub = $doc.Range($pos, $pos + $len)
                $sub.Font.Bold = 1
                $sub.Font.BoldBi = 1
            }
        }
        $pos += $len
    }
}

# Deletes paragraph $paraIndex entirely (text + its own paragraph mark),
# merging what follows up into the previous paragraph's slot.
function Remove-Para($doc, $paraIndex) {
    $doc.Paragraphs($paraIndex).Range.Delete()
}

# Deletes the trailing (now-unwanted) empty paragraph that immediately
# follows paragraph $paraIndex, by removing $paraIndex's own paragraph mark
# (this merges the empty paragraph's mark onto $paraIndex without disturbing
# a possible required final body paragraph mark).
function Remove-ParaMark($doc, $paraIndex) {
    $p = $doc.Paragraphs($paraIndex)
    $e = $p.Range.End
    $doc.Range($e - 1, $e).Delete()
}

# ---------------------------------------------------------------------------
# Resolve every target paragraph index up front (against the pristine,
# unmodified document) using search strings that are unique at this point.
# We will then apply edits from the bottom of the document upward so that
# paragraph-count-changing operations never invalidate an index we still
# need to use.
# ---------------------------------------------------------------------------

$iDataAssessment  = Find-ParaIndex $d "Data Assessment document"
$iCollecting      = Find-ParaIndex $d "Collecting data "
$iSurveyAd        = Find-ParaIndex $d "Survey advertisement through various"
$iCreating        = Find-ParaIndex $d "Creating initial Base DB Structure"
$iInitialDataWork = Find-ParaIndex $d "Initial Data work"
$iRevisiting      = Find-ParaIndex $d "Revisiting the Flow and High-level architecture"

$iGcp             = Find-ParaIndex $d "GCP account and DB structures"
$iConsolidation   = Find-ParaIndex $d "Data Consolidation"
$iDbDesign        = Find-ParaIndex $d "DB Design of level"
$iHighLevel       = Find-ParaIndex $d "High level Flow and design documents"

$iIndicate        = Find-ParaIndex $d "Indicate the"

$iExploring       = Find-ParaIndex $d "Exploring automated pipelines"

$iMidTerm         = Find-ParaIndex $d "Mid Term Presentations"
$iMore            = Find-ParaIndex $d "More "
$iMachine         = Find-ParaIndex $d "Machine learning discussions"

# ---------------------------------------------------------------------------
# Apply edits from the bottom of the document upward.
# ---------------------------------------------------------------------------

# --- Planned Activities for coming week ---

# Remove the trailing empty paragraph at the very end of the document by
# merging $iMachine's own paragraph mark into it.
Remove-ParaMark $d $iMachine

Set-ParaRuns $d $iMachine @(
    @{Text="Initial Dashboard ideas"; Bold=$false}
)

Set-ParaRuns $d $iMore @(
    @{Text="Machine learning "; Bold=$false},
    @{Text="work and Research questions"; Bold=$false}
)

Set-ParaRuns $d $iMidTerm @(
    @{Text="More "; Bold=$false},
    @{Text="Data cleaning logic breakdown"; Bold=$false}
)

# --- Challenges/Disagreements ---
# The trailing "." run is left as-is (unchanged context in the diff).

Set-ParaRunsKeepTail $d $iExploring @(
    @{Text="Exploring automated pipelines for the flow of work"; Bold=$false},
    @{Text=" -some level of API connect is achieved using python between "; Bold=$false},
    @{Text="Colab"; Bold=$false},
    @{Text=" and Big Query"; Bold=$false}
) 1

# --- On Target ---
# The trailing ":" run is left as-is (unchanged context in the diff).

Set-ParaRunsKeepTail $d $iIndicate @(
    @{Text="Indicate the current status of your project"; Bold=$false}
) 1

# --- Specific Output from work ---

Remove-Para $d $iHighLevel

Set-ParaRuns $d $iDbDesign @(
    @{Text="Connections between Big Query and Python "; Bold=$false},
    @{Text="Colab"; Bold=$false},
    @{Text=" (using API) and Tableau"; Bold=$false}
)

Set-ParaRuns $d $iConsolidation @(
    @{Text="Initial Python files and DB scripts"; Bold=$false}
)

Set-ParaRuns $d $iGcp @(
    @{Text="Mid Term presentations"; Bold=$false}
)

# --- Specific Activities that were completed/worked on ---

Remove-Para $d $iRevisiting

Set-ParaRuns $d $iInitialDataWork @(
    @{Text="Initial ML discussions and Research question"; Bold=$false}
)

Set-ParaRuns $d $iCreating @(
    @{Text="Data work – Cleaning, Standardizing and EDA"; Bold=$false}
)

Set-ParaRuns $d $iSurveyAd @(
    @{Text="Extending "; Bold=$false},
    @{Text="initial Base DB Structure in Google Cloud (Big Query)"; Bold=$false},
    @{Text="."; Bold=$false}
)

Set-ParaRuns $d $iCollecting @(
    @{Text="Mid Term Presentations for Capstone"; Bold=$false}
)

Set-ParaRuns $d $iDataAssessment @(
    @{Text="Still "; Bold=$false},
    @{Text="Collecting data "; Bold=$false},
    @{Text="Final Version of Phase 1 Survey"; Bold=$true},
    @{Text=" "; Bold=$false}
)

Write-Output "done"
